$d = $word.ActiveDocument

function Find-ParagraphIndex($doc, $text) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        if ($doc.Paragraphs.Item($i).Range.Text -eq $text) {
            return $i
        }
    }
    return -1
}

# ------------------------------------------------------------------
# 1) Remove the "npm"/"look at page" block that was cut from the
#    workflow (commit message: "removed a few instances of npm").
#    "Npm install", "Npm start", "Look at page", "Stop npm start"
#    (the 4 paragraphs right after "Clone project").
# ------------------------------------------------------------------
$idxNpmInstall = Find-ParagraphIndex $d "Npm install`r"
$idxStopNpmStart = Find-ParagraphIndex $d "Stop npm start`r"
if ($idxNpmInstall -gt 0 -and $idxStopNpmStart -ge $idxNpmInstall) {
    $startPos = $d.Paragraphs.Item($idxNpmInstall).Range.Start
    $endPos = $d.Paragraphs.Item($idxStopNpmStart).Range.End
    $d.Range($startPos, $endPos).Delete()
}

# ------------------------------------------------------------------
# 2) Move "Install git" so it becomes the very first paragraph,
#    immediately before "Git config", with the text split across two
#    runs: "Install git" + a trailing space.
# ------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(1)
$insertRange = $d.Range($p1.Range.Start, $p1.Range.Start)

$installGitXml = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Install git</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertRange.InsertXML($installGitXml)

# Remove the old "Install git" paragraph (now found further down, right
# after the new copy + "Git config").
$idxOldInstallGit = -1
for ($i = 3; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -eq "Install git`r") {
        $idxOldInstallGit = $i
        break
    }
}
if ($idxOldInstallGit -gt 0) {
    $d.Paragraphs.Item($idxOldInstallGit).Range.Delete()
}

# ------------------------------------------------------------------
# 3) Move the "_GoBack" bookmark from the end of the document (after
#    "Merge via merge request into master") to the start of the
#    "Git checkout developer branch" paragraph.
# ------------------------------------------------------------------
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()

$idxDevBranch = Find-ParagraphIndex $d "Git checkout developer branch`r"
if ($idxDevBranch -gt 0) {
    $devStart = $d.Paragraphs.Item($idxDevBranch).Range.Start
    $bmRange = $d.Range($devStart, $devStart)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}
